# Weekly data refresh: insert 3 new price records for "Pepino ensalada"
# (Macroferia Regional de Talca) right before the existing row 710, pushing
# the previously-710..773 block down to 713..776.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 710:773 down by three rows so the three brand-new records can
# be written into the now-empty rows 710:712.
$ws.Rows("710:712").Insert()

# --- New row 710 ---------------------------------------------------------
$ws.Cells.Item(710, 1).Value2 = 5
$ws.Cells.Item(710, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(710, 3).Value2 = "Maule"
$ws.Cells.Item(710, 4).Value2 = 45223
$ws.Cells.Item(710, 5).Value2 = 7
$ws.Cells.Item(710, 6).Value2 = 100112043
$ws.Cells.Item(710, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(710, 8).Value2 = "Sin especificar"
$ws.Cells.Item(710, 9).Value2 = "Primera"
$ws.Cells.Item(710, 10).Value2 = 300
$ws.Cells.Item(710, 11).Value2 = 9000
$ws.Cells.Item(710, 12).Value2 = 9000
$ws.Cells.Item(710, 13).Value2 = 9000
$ws.Cells.Item(710, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(710, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(710, 16).Value2 = 150
$ws.Cells.Item(710, 17).Value2 = 60
$ws.Cells.Item(710, 18).Value2 = "Hortaliza"

# --- New row 711 ---------------------------------------------------------
$ws.Cells.Item(711, 1).Value2 = 5
$ws.Cells.Item(711, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(711, 3).Value2 = "Maule"
$ws.Cells.Item(711, 4).Value2 = 45223
$ws.Cells.Item(711, 5).Value2 = 7
$ws.Cells.Item(711, 6).Value2 = 100112043
$ws.Cells.Item(711, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(711, 8).Value2 = "Sin especificar"
$ws.Cells.Item(711, 9).Value2 = "Primera"
$ws.Cells.Item(711, 10).Value2 = 300
$ws.Cells.Item(711, 11).Value2 = 12000
$ws.Cells.Item(711, 12).Value2 = 12000
$ws.Cells.Item(711, 13).Value2 = 12000
$ws.Cells.Item(711, 14).Value2 = "`$/caja 80 unidades"
$ws.Cells.Item(711, 15).Value2 = "Región del Maule"
$ws.Cells.Item(711, 16).Value2 = 150
$ws.Cells.Item(711, 17).Value2 = 80
$ws.Cells.Item(711, 18).Value2 = "Hortaliza"

# --- New row 712 ---------------------------------------------------------
$ws.Cells.Item(712, 1).Value2 = 5
$ws.Cells.Item(712, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(712, 3).Value2 = "Maule"
$ws.Cells.Item(712, 4).Value2 = 45223
$ws.Cells.Item(712, 5).Value2 = 7
$ws.Cells.Item(712, 6).Value2 = 100112043
$ws.Cells.Item(712, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(712, 8).Value2 = "Sin especificar"
$ws.Cells.Item(712, 9).Value2 = "Segunda"
$ws.Cells.Item(712, 10).Value2 = 200
$ws.Cells.Item(712, 11).Value2 = 8000
$ws.Cells.Item(712, 12).Value2 = 8000
$ws.Cells.Item(712, 13).Value2 = 8000
$ws.Cells.Item(712, 14).Value2 = "`$/caja 100 unidades"
$ws.Cells.Item(712, 15).Value2 = "Región del Maule"
$ws.Cells.Item(712, 16).Value2 = 80
$ws.Cells.Item(712, 17).Value2 = 100
$ws.Cells.Item(712, 18).Value2 = "Hortaliza"
